$wb = $excel.ActiveWorkbook

# ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 46784.2
$ws.Range("J3").Value = 46784.2
$ws.Range("L3").Value = 46784.2
$ws.Range("N3").Value = -47012.2
$ws.Range("H20").Value = 632.25
$ws.Range("I20").Value = 632.25
$ws.Range("K20").Value = 632.25
$ws.Range("M20").Value = -402.25
$ws.Range("H33").Value = 1102.1333
$ws.Range("I33").Value = 181.875
$ws.Range("J33").Value = 2153.8572
$ws.Range("K33").Value = 181.875
$ws.Range("L33").Value = 2153.8572
$ws.Range("M33").Value = 47.125
$ws.Range("N33").Value = -2611.8572
$ws.Range("H35").Value = 632.25
$ws.Range("I35").Value = 632.25
$ws.Range("K35").Value = 632.25
$ws.Range("M35").Value = -253.25
$ws.Range("H47").Value = 11111
$ws.Range("I47").Value = 11111
$ws.Range("K47").Value = 11111
$ws.Range("M47").Value = -10139
$ws.Range("H101").Value = 1896.2
$ws.Range("I101").Value = 440.22223
$ws.Range("J101").Value = 15000
$ws.Range("K101").Value = 1320.66669
$ws.Range("L101").Value = 45000
$ws.Range("M101").Value = 301.33331
$ws.Range("N101").Value = -48244
$ws.Range("H102").Value = 46784.2
$ws.Range("J102").Value = 46784.2
$ws.Range("L102").Value = 46784.2
$ws.Range("N102").Value = -53274.2
$ws.Range("H138").Value = 3568.7637
$ws.Range("I138").Value = 1349.4
$ws.Range("J138").Value = 4061.9556
$ws.Range("K138").Value = 4048.2
$ws.Range("L138").Value = 12185.8668
$ws.Range("M138").Value = 1091.8
$ws.Range("N138").Value = -22465.8668

# ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1456.25
$ws.Range("I26").Value = 950
$ws.Range("J26").Value = 5000
$ws.Range("K26").Value = 950
$ws.Range("L26").Value = 5000
$ws.Range("M26").Value = -620
$ws.Range("N26").Value = -5660
$ws.Range("H32").Value = 3033.8723
$ws.Range("I32").Value = 2946.6978
$ws.Range("J32").Value = 3971
$ws.Range("K32").Value = 2946.6978
$ws.Range("L32").Value = 3971
$ws.Range("M32").Value = -2659.6978
$ws.Range("N32").Value = -4545
$ws.Range("H45").Value = 109092620
$ws.Range("J45").Value = 111113000
$ws.Range("L45").Value = 111113000
$ws.Range("N45").Value = -111113754
$ws.Range("H74").Value = 1839.7222
$ws.Range("I74").Value = 1932.7693
$ws.Range("J74").Value = 1597.8
$ws.Range("K74").Value = 1932.7693
$ws.Range("L74").Value = 1597.8
$ws.Range("M74").Value = -1058.7693
$ws.Range("N74").Value = -3345.8
$ws.Range("H77").Value = 1839.7222
$ws.Range("I77").Value = 1932.7693
$ws.Range("J77").Value = 1597.8
$ws.Range("K77").Value = 9663.8465
$ws.Range("L77").Value = 7989
$ws.Range("M77").Value = -5295.8465
$ws.Range("N77").Value = -16725
$ws.Range("H97").Value = 1497.6
$ws.Range("I97").Value = 1537.2142
$ws.Range("J97").Value = 1405.1666
$ws.Range("K97").Value = 1537.2142
$ws.Range("L97").Value = 1405.1666
$ws.Range("M97").Value = -1041.2142
$ws.Range("N97").Value = -2397.1666

# BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2662.4546
$ws.Range("I94").Value = 2700.2273
$ws.Range("J94").Value = 2586.9092
$ws.Range("K94").Value = 2700.2273
$ws.Range("L94").Value = 2586.9092
$ws.Range("M94").Value = -2249.2273
$ws.Range("N94").Value = -3488.9092
$ws.Range("H120").Value = 71280
$ws.Range("J120").Value = 71280
$ws.Range("L120").Value = 71280
$ws.Range("N120").Value = -80956

# CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3452.5
$ws.Range("I16").Value = 3452.5
$ws.Range("K16").Value = 3452.5
$ws.Range("M16").Value = -3165.5
$ws.Range("H62").Value = 8546.200000000001
$ws.Range("I62").Value = 7465.8887
$ws.Range("K62").Value = 7465.8887
$ws.Range("M62").Value = -6841.8887
$ws.Range("H65").Value = 8546.200000000001
$ws.Range("I65").Value = 7465.8887
$ws.Range("K65").Value = 37329.4435
$ws.Range("M65").Value = -34209.4435
$ws.Range("H99").Value = 2897.2222
$ws.Range("I99").Value = 2867.8572
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2867.8572
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -1369.8572
$ws.Range("N99").Value = -5996
$ws.Range("H113").Value = 3452.5
$ws.Range("I113").Value = 3452.5
$ws.Range("K113").Value = 3452.5
$ws.Range("M113").Value = -1282.5
$ws.Range("H126").Value = 2897.2222
$ws.Range("I126").Value = 2867.8572
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8603.571599999999
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -6133.571599999999
$ws.Range("N126").Value = -13940
$ws.Range("H134").Value = 3822.9333
$ws.Range("I134").Value = 4005.7693
$ws.Range("K134").Value = 12017.3079
$ws.Range("M134").Value = -9482.3079

# CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 13944.444
$ws.Range("J88").Value = 13666.667
$ws.Range("L88").Value = 41000.001
$ws.Range("N88").Value = -41856.001
$ws.Range("H91").Value = 13944.444
$ws.Range("J91").Value = 13666.667
$ws.Range("L91").Value = 41000.001
$ws.Range("N91").Value = -43964.001

# GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 33679
$ws.Range("J42").Value = 54995
$ws.Range("L42").Value = 54995
$ws.Range("N42").Value = -55965
$ws.Range("H115").Value = 33679
$ws.Range("J115").Value = 54995
$ws.Range("L115").Value = 54995
$ws.Range("N115").Value = -57345

# LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5599.3335
$ws.Range("I7").Value = 3999.25
$ws.Range("J7").Value = 8799.5
$ws.Range("K7").Value = 3999.25
$ws.Range("L7").Value = 8799.5
$ws.Range("M7").Value = -3887.25
$ws.Range("N7").Value = -9023.5
$ws.Range("H40").Value = 6725.1787
$ws.Range("J40").Value = 7933.1333
$ws.Range("L40").Value = 7933.1333
$ws.Range("N40").Value = -8205.133300000001
$ws.Range("H122").Value = 4950
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 5599.3335
$ws.Range("I126").Value = 3999.25
$ws.Range("J126").Value = 8799.5
$ws.Range("K126").Value = 11997.75
$ws.Range("L126").Value = 26398.5
$ws.Range("M126").Value = -9527.75
$ws.Range("N126").Value = -31338.5
$ws.Range("H132").Value = 4588.1113
$ws.Range("I132").Value = 4549.125
$ws.Range("J132").Value = 4900
$ws.Range("K132").Value = 13647.375
$ws.Range("L132").Value = 14700
$ws.Range("M132").Value = -11117.375
$ws.Range("N132").Value = -19760
$ws.Range("H136").Value = 4957.6523
$ws.Range("I136").Value = 3639.7368
$ws.Range("J136").Value = 5885.074
$ws.Range("K136").Value = 10919.2104
$ws.Range("L136").Value = 17655.222
$ws.Range("M136").Value = -8369.2104
$ws.Range("N136").Value = -22755.222
$ws.Range("H138").Value = 64449
$ws.Range("J138").Value = 64449
$ws.Range("L138").Value = 64449
$ws.Range("N138").Value = -74729
$ws.Range("H141").Value = 52500
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 75000
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 75000
$ws.Range("M141").Value = -24820
$ws.Range("N141").Value = -85360

# WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H116").Value = 221998.33
$ws.Range("J116").Value = 221998.33
$ws.Range("L116").Value = 221998.33
$ws.Range("N116").Value = -231176.33
$ws.Range("H126").Value = 1877.8572
$ws.Range("I126").Value = 1857.6666
$ws.Range("K126").Value = 5572.9998
$ws.Range("M126").Value = -3102.9998
$ws.Range("H139").Value = 94999
$ws.Range("I139").Value = 90000
$ws.Range("J139").Value = 99998
$ws.Range("K139").Value = 90000
$ws.Range("L139").Value = 99998
$ws.Range("M139").Value = -84860
$ws.Range("N139").Value = -110278
